$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E; this shifts the old E..AD data right to F..AE
# and carries the D-column's "s=1" (scientific-notation) style into the new
# E cells for the summary/process rows, which is exactly where the new
# citation column needs it.
$ws.Columns("E").Insert()

# The FDR threshold text used on all the "child" gene rows changed from
# ">0.05" to ">0.1".
$ws.Range("D4").Value = ">0.1"
$ws.Range("D5").Value = ">0.1"
$ws.Range("D6").Value = ">0.1"
$ws.Range("D7").Value = ">0.1"
$ws.Range("D8").Value = ">0.1"
$ws.Range("D9").Value = ">0.1"
$ws.Range("D10").Value = ">0.1"
$ws.Range("D11").Value = ">0.1"
$ws.Range("D12").Value = ">0.1"
$ws.Range("D13").Value = ">0.1"
$ws.Range("D14").Value = ">0.1"
$ws.Range("D15").Value = ">0.1"
$ws.Range("D17").Value = ">0.1"
$ws.Range("D18").Value = ">0.1"
$ws.Range("D19").Value = ">0.1"
$ws.Range("D20").Value = ">0.1"
$ws.Range("D21").Value = ">0.1"
$ws.Range("D22").Value = ">0.1"
$ws.Range("D23").Value = ">0.1"
$ws.Range("D24").Value = ">0.1"
$ws.Range("D26").Value = ">0.1"
$ws.Range("D27").Value = ">0.1"
$ws.Range("D28").Value = ">0.1"
$ws.Range("D30").Value = ">0.1"
$ws.Range("D31").Value = ">0.1"
$ws.Range("D32").Value = ">0.1"
$ws.Range("D33").Value = ">0.1"
$ws.Range("D35").Value = ">0.1"
$ws.Range("D36").Value = ">0.1"
$ws.Range("D37").Value = ">0.1"
$ws.Range("D38").Value = ">0.1"
$ws.Range("D39").Value = ">0.1"
$ws.Range("D41").Value = ">0.1"
$ws.Range("D42").Value = ">0.1"
$ws.Range("D43").Value = ">0.1"
$ws.Range("D44").Value = ">0.1"

# New column E holds a DOI / citation reference for each "process" summary
# row (the rows that already carried a Process name in column C).
$ws.Range("E2").Value = "10.1091/mbc.e09-12-1031"
$ws.Range("E3").Value = "10.1073/pnas.82.21.7193"
$ws.Range("E16").Value = "(none)"
$ws.Range("E25").Value = "(none)"
$ws.Range("E29").Value = "10.1016/j.semcdb.2006.10.011"
$ws.Range("E34").Value = "(none)"
$ws.Range("E40").Value = "10.1080/02713680500477347"
$ws.Range("E45").Value = "10.1042/BJ20040347"
$ws.Range("E46").Value = "10.1016/0039-6257(88)90095-1"

# Update the view: scroll back to the top and move the active selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E49").Select()
